$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Shrink the existing bold/italic run from the old filename text
#    down to just the ".py" suffix, leaving its bold/italic rPr intact.
#    The old filename text is unique in the document, so a single
#    (ReplaceOne) Find/Replace is safe and precise.
# ------------------------------------------------------------------
$oldName = "obtainDataINITandPositioning.py"
$find1 = $d.Content
$ok1 = $find1.Find.Execute($oldName, $true, $false, $false, $false, $false, $true, 1, $false, ".py", 1)
if (-not $ok1) {
    throw "Could not find the original run text '$oldName' to shrink to '.py'."
}

# ------------------------------------------------------------------
# 2. Re-locate the now-shrunk run (unique text ".py" immediately
#    preceded by "from ") and insert a brand-new run of text right
#    before it: "datasetINIT23_main". InsertBefore merges the new
#    characters into the immediately-preceding run while they share
#    identical formatting, so we immediately re-select exactly the
#    newly inserted characters and give them their own distinct
#    character formatting (Arial font + dark grey colour) which forces
#    the serializer to split them into their own <w:r>.
# ------------------------------------------------------------------
$find2 = $d.Content
$ok2 = $find2.Find.Execute(".py", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok2) {
    throw "Could not re-locate the shrunk '.py' run."
}

$insertPos = $find2.Start
$insertPoint = $d.Range($insertPos, $insertPos)
$newText = "datasetINIT23_main"
$insertPoint.InsertBefore($newText)

$newRun = $d.Range($insertPos, $insertPos + $newText.Length)
$newRun.Font.Name = "Arial"
$newRun.Font.Color = 2236962   # RGB 0x222222 (34,34,34) as a Word BGR long

Write-Output ("Inserted run text: " + $newRun.Text)
Write-Output ("Trailing run text: " + $d.Range($insertPos + $newText.Length, $insertPos + $newText.Length + 3).Text)
